$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 393; this shifts old rows 393..519 down to 394..520,
# carrying their formatting (incl. the date-style column D) with them.
$ws.Rows("393:393").Insert()

# Populate the newly inserted row 393 with the new data point.
# Static columns mirror the surrounding rows (now 392 / 394).
$ws.Range("A393").Value = 3
$ws.Range("B393").Value = "Femacal de La Calera"
$ws.Range("C393").Value = "Coquimbo"
$ws.Range("D393").Value = 44985
$ws.Range("D393").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E393").Value = 5
$ws.Range("F393").Value = 100112040
$ws.Range("G393").Value = "Cilantro"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 80
$ws.Range("K393").Value = 5500
$ws.Range("L393").Value = 5500
$ws.Range("M393").Value = 5500
$ws.Range("N393").Value = "$/docena de atados (3 kilos)"
$ws.Range("O393").Value = "Provincia de Quillota"
$ws.Range("P393").Value = 1833
$ws.Range("Q393").Value = 3
$ws.Range("R393").Value = "Hortaliza"
